$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 54
$dstRow = 55

# Copy formatting (styles) from the last existing data row to the new row,
# matching the s="1" border/bold style on column A and s="2" date format on column E.
$src = $ws.Range("A" + $srcRow + ":V" + $srcRow)
$dst = $ws.Range("A" + $dstRow + ":V" + $dstRow)
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item($dstRow, 1).Value = 54
$ws.Cells.Item($dstRow, 2).Value = "belgium"
$ws.Cells.Item($dstRow, 3).Value = "jupiler-pro-league"
$ws.Cells.Item($dstRow, 4).Value = "2023-2024"
$ws.Cells.Item($dstRow, 5).Value = 45191.86458333334
$ws.Cells.Item($dstRow, 6).Value = "St. Liege"
$ws.Cells.Item($dstRow, 7).Value = 0
$ws.Cells.Item($dstRow, 8).Value = "Westerlo"
$ws.Cells.Item($dstRow, 9).Value = 0
$ws.Cells.Item($dstRow, 10).Value = 1.99
$ws.Cells.Item($dstRow, 11).Value = "17/09/2023 12:42"
$ws.Cells.Item($dstRow, 12).Value = 2.14
$ws.Cells.Item($dstRow, 13).Value = "22/09/2023 20:44"
$ws.Cells.Item($dstRow, 14).Value = 3.66
$ws.Cells.Item($dstRow, 15).Value = "17/09/2023 12:42"
$ws.Cells.Item($dstRow, 16).Value = 3.76
$ws.Cells.Item($dstRow, 17).Value = "22/09/2023 20:44"
$ws.Cells.Item($dstRow, 18).Value = 3.72
$ws.Cells.Item($dstRow, 19).Value = "17/09/2023 12:42"
$ws.Cells.Item($dstRow, 20).Value = 3.32
$ws.Cells.Item($dstRow, 21).Value = "22/09/2023 20:44"
$ws.Cells.Item($dstRow, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/st-liege-westerlo/MFrCP7qC/"
